# "service add square method" -- add a computed "Square" (width * height)
# column derived from the existing "Resolution" column (e.g. "1900x475"),
# plus two free-text note rows below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Free-text notes appended below the table (written first so the shared
#     -string table order matches the authored workbook) ---
$ws.Range("B51").Value = "это назвать 1-2"
$ws.Range("A53").Value = "остальное нормировать"

# --- New column header: G1 = "Square" (bold, matches the other header cells) ---
$ws.Range("G1").Value = "Square"
$ws.Range("G1").Font.Bold = $true

# --- Compute Square = width * height from the "Resolution" column (D) for
#     every data row and write it into the new "Square" column (G) ---
for ($r = 2; $r -le 50; $r++) {
    $resolution = $ws.Cells.Item($r, 4).Value2
    $parts = $resolution.Split("x")
    $width = [double]$parts[0]
    $height = [double]$parts[1]
    $ws.Cells.Item($r, 7).Value = $width * $height
}

# --- View bookkeeping (selection + zoom) to mirror the author's session ---
$ws.Range("L6").Select()
$excel.ActiveWindow.Zoom = 94
